# Update betting-odds/score values for rows 5-12 on Sheet1
# per updated FlashScore data (2024-11-05).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5
$ws.Range("G5").Value = 3.1
$ws.Range("I5").Value = 2.25
$ws.Range("J5").Value = 3.6
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.85
$ws.Range("X5").Value = 15
$ws.Range("Y5").Value = 11
$ws.Range("AC5").Value = 10
$ws.Range("AL5").Value = 19
$ws.Range("BC5").Value = 126
# Row 6
$ws.Range("G6").Value = 2.9
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 2.5
$ws.Range("J6").Value = 3.5
$ws.Range("L6").Value = 3.1
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 9
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 3.4
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.73
$ws.Range("S6").Value = 1.44
$ws.Range("T6").Value = 2.63
$ws.Range("U6").Value = 1.8
$ws.Range("V6").Value = 1.91
$ws.Range("W6").Value = 9
$ws.Range("Z6").Value = 29
$ws.Range("AC6").Value = 9
$ws.Range("AF6").Value = 51
$ws.Range("AG6").Value = 251
$ws.Range("AI6").Value = 12
$ws.Range("AJ6").Value = 10
$ws.Range("AK6").Value = 23
$ws.Range("AL6").Value = 21
$ws.Range("AN6").Value = 4.75
$ws.Range("AS6").Value = 201
$ws.Range("AT6").Value = 2.63
$ws.Range("AW6").Value = 4.5
$ws.Range("BC6").Value = 126
# Row 7
$ws.Range("I7").Value = 1.83
$ws.Range("AL7").Value = 13
# Row 8
$ws.Range("G8").Value = 2
$ws.Range("I8").Value = 3.6
$ws.Range("K8").Value = 2.2
$ws.Range("Q8").Value = 2.03
$ws.Range("R8").Value = 1.83
$ws.Range("S8").Value = 1.37
$ws.Range("U8").Value = 1.8
$ws.Range("V8").Value = 1.91
$ws.Range("X8").Value = 9.5
$ws.Range("AE8").Value = 15
$ws.Range("AG8").Value = 251
$ws.Range("AH8").Value = 11
$ws.Range("AI8").Value = 19
$ws.Range("AJ8").Value = 13
$ws.Range("AO8").Value = 11
$ws.Range("AX8").Value = 21
$ws.Range("AY8").Value = 29
# Row 9
$ws.Range("G9").Value = 2.35
$ws.Range("I9").Value = 2.8
$ws.Range("J9").Value = 3
$ws.Range("K9").Value = 2.25
$ws.Range("L9").Value = 3.25
$ws.Range("S9").Value = 1.3
$ws.Range("AO9").Value = 13
$ws.Range("AP9").Value = 21
$ws.Range("AU9").Value = 7
# Row 10
$ws.Range("O10").Value = 1.22
$ws.Range("P10").Value = 4.33
$ws.Range("Q10").Value = 1.73
$ws.Range("R10").Value = 2.1
$ws.Range("S10").Value = 1.3
# Row 11
$ws.Range("G11").Value = 1.95
$ws.Range("I11").Value = 3.5
$ws.Range("J11").Value = 2.6
$ws.Range("L11").Value = 4
$ws.Range("S11").Value = 1.33
$ws.Range("Y11").Value = 8.5
$ws.Range("Z11").Value = 17
$ws.Range("AM11").Value = 34
$ws.Range("AN11").Value = 4
$ws.Range("AY11").Value = 26
$ws.Range("BA11").Value = 81
# Row 12
$ws.Range("I12").Value = 2.8
$ws.Range("M12").Value = 1.05
$ws.Range("N12").Value = 11
$ws.Range("X12").Value = 13
$ws.Range("Y12").Value = 10
$ws.Range("AH12").Value = 9.5
$ws.Range("AL12").Value = 21
$ws.Range("AP12").Value = 23
